$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.946.95"
$ws.Range("D2").Style = $ws.Range("A1").Style
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.448.63"
$ws.Range("D3").Style = $ws.Range("A1").Style
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "512.30"
$ws.Range("D5").Style = $ws.Range("A1").Style
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.14"
$ws.Range("D6").Style = $ws.Range("A1").Style
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("D8").Style = $ws.Range("A1").Style
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.466.58"
$ws.Range("D9").Style = $ws.Range("A1").Style
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0962"
$ws.Range("D10").Style = $ws.Range("A1").Style
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.156"
$ws.Range("D11").Style = $ws.Range("A1").Style
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.21"
$ws.Range("D12").Style = $ws.Range("A1").Style
$ws.Range("E12").Value = "  -2.69%  "
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.887.78"
$ws.Range("D14").Style = $ws.Range("A1").Style
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.900.79"
$ws.Range("D15").Style = $ws.Range("A1").Style
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.00"
$ws.Range("D16").Style = $ws.Range("A1").Style
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.460.57"
$ws.Range("D18").Style = $ws.Range("A1").Style
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.55"
$ws.Range("D19").Style = $ws.Range("A1").Style
$ws.Range("E19").Value = "  -2.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "318.84"
$ws.Range("D20").Style = $ws.Range("A1").Style
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("E21").Value = "  -1.43%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.97"
$ws.Range("D23").Style = $ws.Range("A1").Style
$ws.Range("E23").Value = "  +3.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.42"
$ws.Range("D24").Style = $ws.Range("A1").Style
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.404"
$ws.Range("D25").Style = $ws.Range("A1").Style
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.160"
$ws.Range("D26").Style = $ws.Range("A1").Style
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.976"
$ws.Range("D27").Style = $ws.Range("A1").Style
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.97"
$ws.Range("D29").Style = $ws.Range("A1").Style
$ws.Range("E29").Value = "  +0.97%  "
$ws.Range("E30").Value = "  -3.00%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.67"
$ws.Range("D31").Style = $ws.Range("A1").Style
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.17"
$ws.Range("D32").Style = $ws.Range("A1").Style
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.19"
$ws.Range("D33").Style = $ws.Range("A1").Style
$ws.Range("E33").Value = "  -2.43%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.995"
$ws.Range("D35").Style = $ws.Range("A1").Style
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.27"
$ws.Range("D37").Style = $ws.Range("A1").Style
$ws.Range("E37").Value = "  -3.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.92"
$ws.Range("D38").Style = $ws.Range("A1").Style
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.56"
$ws.Range("D39").Style = $ws.Range("A1").Style
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.762"
$ws.Range("D41").Style = $ws.Range("A1").Style
$ws.Range("E41").Value = "  -4.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "272.56"
$ws.Range("D42").Style = $ws.Range("A1").Style
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.39"
$ws.Range("D43").Style = $ws.Range("A1").Style
$ws.Range("E43").Value = "  -2.88%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.00"
$ws.Range("D44").Style = $ws.Range("A1").Style
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.587"
$ws.Range("D45").Style = $ws.Range("A1").Style
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.68"
$ws.Range("D47").Style = $ws.Range("A1").Style
$ws.Range("E47").Value = "  -5.26%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.31"
$ws.Range("D49").Style = $ws.Range("A1").Style
$ws.Range("E49").Value = "  -3.37%  "
$ws.Range("E50").Value = "  -2.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.71"
$ws.Range("D51").Style = $ws.Range("A1").Style
$ws.Range("E51").Value = "  -2.84%  "
